$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 <- old row 3 values
$ws.Range("D2").Value = 44229
$ws.Range("K2").Value = "Santina"
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 250
$ws.Range("N2").Value = 6500
$ws.Range("O2").Value = 7000
$ws.Range("P2").Value = 6750
$ws.Range("Q2").Value = "`$/bandeja 5 kilos"
$ws.Range("R2").Value = "Provincia de Curicó"
$ws.Range("S2").Value = 1350
$ws.Range("T2").Value = 5

# Row 3 <- old row 5 values
$ws.Range("D3").Value = 44161
$ws.Range("K3").Value = "Bing"
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 160
$ws.Range("N3").Value = 39000
$ws.Range("O3").Value = 40000
$ws.Range("P3").Value = 39500
$ws.Range("Q3").Value = "`$/caja 20 kilos"
$ws.Range("R3").Value = "Provincia de Curicó"
$ws.Range("S3").Value = 1975
$ws.Range("T3").Value = 20

# Row 4 <- old row 2 values
$ws.Range("D4").Value = 44208
$ws.Range("K4").Value = "Lapins"
$ws.Range("L4").Value = "Segunda"
$ws.Range("M4").Value = 200
$ws.Range("N4").Value = 10500
$ws.Range("O4").Value = 11000
$ws.Range("P4").Value = 10750
$ws.Range("Q4").Value = "`$/bandeja 12 kilos"
$ws.Range("R4").Value = "Provincia de Curicó"
$ws.Range("S4").Value = 896
$ws.Range("T4").Value = 12

# Row 5 <- old row 4 values
$ws.Range("D5").Value = 44210
$ws.Range("K5").Value = "Rainier"
$ws.Range("L5").Value = "Segunda"
$ws.Range("M5").Value = 250
$ws.Range("N5").Value = 21000
$ws.Range("O5").Value = 22000
$ws.Range("P5").Value = 21500
$ws.Range("Q5").Value = "`$/caja 18 kilos"
$ws.Range("R5").Value = "Región de O'Higgins"
$ws.Range("S5").Value = 1194
$ws.Range("T5").Value = 18
